$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (23-37) appended to the ExtractedScans table, columns A-I.
# Column order: RecNumber, Label, Fat, Lean, FreeWater, TotalWater, TimeDateDura, Accumulation, Weight
$rows = @(
    @(1,  9401, 3.78,  35.08, 0, 0, "09:13:28 Dec 18, 2025; 27; ems", 1, 39),
    @(2,  9403, 12.41, 33.44, 0, 0, "09:14:35 Dec 18, 2025; 31; ems", 1, 45.7),
    @(3,  9411, 5.41,  30.82, 0, 0, "09:15:44 Dec 18, 2025; 31; ems", 1, 37.1),
    @(4,  9407, 15.81, 38.47, 0, 0, "09:17:04 Dec 18, 2025; 27; ems", 1, 53.9),
    @(5,  9418, 6.44,  36.48, 0, 0, "09:18:05 Dec 18, 2025; 31; ems", 1, 44.3),
    @(6,  9410, 19.16, 43.26, 0, 0, "09:19:27 Dec 18, 2025; 31; ems", 1, 62.3),
    @(7,  9414, 3.43,  24.48, 0, 0, "09:20:59 Dec 18, 2025; 31; ems", 1, 28.6),
    @(8,  9417, 9.14,  35.69, 0, 0, "09:22:05 Dec 18, 2025; 32; ems", 1, 45),
    @(9,  9405, 5.08,  24.8,  0, 0, "09:23:12 Dec 18, 2025; 27; ems", 1, 30.5),
    @(10, 9402, 9.59,  27.34, 0, 0, "09:24:09 Dec 18, 2025; 27; ems", 1, 37.2),
    @(11, 9404, 1.14,  25.07, 0, 0, "09:25:28 Dec 18, 2025; 36; ems", 1, 27.5),
    @(12, 9406, 0.82,  25.99, 0, 0, "09:26:43 Dec 18, 2025; 31; ems", 1, 27.8),
    @(13, 9400, 12.06, 38.2,  0, 0, "09:27:58 Dec 18, 2025; 27; ems", 1, 50.6),
    @(14, 9412, 2.27,  31.12, 0, 0, "09:29:12 Dec 18, 2025; 31; ems", 1, 34.7),
    @(15, 9415, 12.88, 36.66, 0, 0, "09:30:23 Dec 18, 2025; 27; ems", 1, 49.8)
)

$startRow = 23
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}

# New font/style (Helvetica 10) applied across the newly added block, matching
# the workbook's new cellXfs entry (fontId=1, applyFont=1).
$endRow = $startRow + $rows.Count - 1
$newRange = $ws.Range("A$startRow`:I$endRow")
$newRange.Font.Name = "Helvetica"
$newRange.Font.Size = 10

# Column G widened slightly to accommodate the new content.
$ws.Columns("G").ColumnWidth = 29

# Move the active selection to just past the newly appended data, as in the
# saved workbook (G38).
$ws.Range("G38").Select() | Out-Null
